$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("W1").Value = 0.82771046300565021
$ws.Range("BP1").Value = 0.78402289885913246
$ws.Range("C2").Value = 0.94862301804875837
$ws.Range("A3").Value = 0.88860284244794352
$ws.Range("D3").Value = 0.70727105912414845
$ws.Range("BG3").Value = 0.93888149240445529
$ws.Range("B4").Value = 0.87834147519180394
$ws.Range("E4").Value = 0.8878175166418798
$ws.Range("F4").Value = 0.79043653901171362
$ws.Range("G6").Value = 0.69355427985150297
$ws.Range("H7").Value = 0.78025980645017168
$ws.Range("I7").Value = 0.77607356003902939
$ws.Range("F8").Value = 0.99578426808316411
$ws.Range("J8").Value = 0.87302430950926002
$ws.Range("T8").Value = 0.94603207668818801
$ws.Range("X9").Value = 0.55666378869579236
$ws.Range("BF9").Value = 0.8117009502960516
$ws.Range("K10").Value = 0.97085656003627974
$ws.Range("AJ10").Value = 0.73185105773155157
$ws.Range("AL10").Value = 0.7864308655702007
$ws.Range("L11").Value = 0.86687549935416452
$ws.Range("BM11").Value = 0.95518912027042102
$ws.Range("AB12").Value = 0.83415773400792026
$ws.Range("L13").Value = 0.99450592790994197
$ws.Range("N13").Value = 0.93456938075925
$ws.Range("U14").Value = 0.99446474589436495
$ws.Range("E15").Value = 0.61988334025775826
$ws.Range("BK16").Value = 0.95929348526596103
$ws.Range("O17").Value = 0.96288137994885281
$ws.Range("S17").Value = 0.87595378978444804
$ws.Range("BF17").Value = 0.91830806950810806
$ws.Range("P18").Value = 0.82040032304159449
$ws.Range("T18").Value = 0.83672367251756863
$ws.Range("R19").Value = 0.93820865905058481
$ws.Range("X19").Value = 0.64706268010145962
$ws.Range("AK20").Value = 0.6334548474597631
$ws.Range("AU20").Value = 0.69841858531996048
$ws.Range("W21").Value = 0.93860123011058949
$ws.Range("AD21").Value = 0.98582590224283873
$ws.Range("N22").Value = 0.85838296458716434
$ws.Range("Q22").Value = 0.90761333917302545
$ws.Range("U22").Value = 0.98091458764897377
$ws.Range("V24").Value = 0.87404933101003368
$ws.Range("X25").Value = 0.72287870547881439
$ws.Range("AO25").Value = 0.96832772963870495
$ws.Range("BA25").Value = 0.91954970987354767
$ws.Range("AK26").Value = 0.70618368143086507
$ws.Range("BO26").Value = 0.62010551326870655
$ws.Range("Z27").Value = 0.58232724975435934
$ws.Range("AB27").Value = 0.82840260040252856
$ws.Range("AC27").Value = 0.6698266092613181
$ws.Range("AQ27").Value = 0.81324791285960429
$ws.Range("T29").Value = 0.82680782718940793
$ws.Range("AB29").Value = 0.95622611532935597
$ws.Range("AD29").Value = 0.93965105060137144
$ws.Range("AZ29").Value = 0.87635832580788597
$ws.Range("Q30").Value = 0.88720466491378591
$ws.Range("AB30").Value = 0.88055437903827716
$ws.Range("AF30").Value = 0.64651639547966977
$ws.Range("N31").Value = 0.72635159214604361
$ws.Range("AF31").Value = 0.88167754280267796
$ws.Range("BD31").Value = 0.7547850225177033
$ws.Range("AH33").Value = 0.9545769265981161
$ws.Range("AI33").Value = 0.81391308358699621
$ws.Range("AF34").Value = 0.92716212418453781
$ws.Range("AJ34").Value = 0.86710842872751326
$ws.Range("AP34").Value = 0.82058904160793489
$ws.Range("Y35").Value = 0.87681075734276359
$ws.Range("AH35").Value = 0.77974266205967535
$ws.Range("AV35").Value = 0.86133595734907809
$ws.Range("AL36").Value = 0.79529091413788211
$ws.Range("AK38").Value = 0.99638645608881715
$ws.Range("AL39").Value = 0.81164160252050477
$ws.Range("AO39").Value = 0.87118616252546932
$ws.Range("I40").Value = 0.79963845693586322
$ws.Range("M40").Value = 0.97227883251682923
$ws.Range("AL40").Value = 0.55577071582208326
$ws.Range("AO40").Value = 0.96666107059148954
$ws.Range("BH40").Value = 0.7492501949171908
$ws.Range("U42").Value = 0.67373615332270864
$ws.Range("AN42").Value = 0.75834318602829009
$ws.Range("AO43").Value = 0.61245594075328447
$ws.Range("AS43").Value = 0.96996286087124095
$ws.Range("L44").Value = 0.90420178630829318
$ws.Range("AM44").Value = 0.79191944144755344
$ws.Range("AR45").Value = 0.54706223524962627
$ws.Range("AT45").Value = 0.87571518691290917
$ws.Range("O46").Value = 0.85939640515763149
$ws.Range("AU46").Value = 0.78452087394168435
$ws.Range("AP47").Value = 0.96861615353420283
$ws.Range("AS47").Value = 0.94753660178093324
$ws.Range("S48").Value = 0.92824516698617465
$ws.Range("AT48").Value = 0.65082106648748894
$ws.Range("AW48").Value = 0.90351227634970077
$ws.Range("AF49").Value = 0.92618056089341794
$ws.Range("BM49").Value = 0.83216511161051987
$ws.Range("AN50").Value = 0.68240537739462059
$ws.Range("AV50").Value = 0.82618273202399217
$ws.Range("AW51").Value = 0.99333483864424255
$ws.Range("AZ51").Value = 0.97958823671840012
$ws.Range("BA51").Value = 0.6895357877997379
$ws.Range("AX52").Value = 0.67392286978954208
$ws.Range("BA52").Value = 0.98909521286707613
$ws.Range("W53").Value = 0.87646387616837029
$ws.Range("BA54").Value = 0.93592944701374248
$ws.Range("BD54").Value = 0.99462594643578961
$ws.Range("AC55").Value = 0.73327810070702637
$ws.Range("AO55").Value = 0.61396315207140972
$ws.Range("BA55").Value = 0.70657927206555149
$ws.Range("BB55").Value = 0.63123312692347011
$ws.Range("BE55").Value = 0.78195465922366403
$ws.Range("BC56").Value = 0.96065414846499242
$ws.Range("BF56").Value = 0.60106288256758866
$ws.Range("W58").Value = 0.88347364094977543
$ws.Range("BE58").Value = 0.63044937508066079
$ws.Range("BH58").Value = 0.97656628660450884
$ws.Range("BJ58").Value = 0.60238646206859658
$ws.Range("BE59").Value = 0.7741057508176814
$ws.Range("BG60").Value = 0.82523493043229612
$ws.Range("BI60").Value = 0.95394573830354001
$ws.Range("Y61").Value = 0.76192202695470712
$ws.Range("BG61").Value = 0.9886607632647576
$ws.Range("BJ61").Value = 0.88334831761646471
$ws.Range("AL62").Value = 0.97664899238020642
$ws.Range("BL63").Value = 0.6952853667946759
$ws.Range("BM63").Value = 0.96681870033230177
$ws.Range("K64").Value = 0.93304909011324644
$ws.Range("BM64").Value = 0.8710661897710803
$ws.Range("D65").Value = 0.94830801128409381
$ws.Range("E66").Value = 0.9784710397849995
$ws.Range("AM66").Value = 0.96699844204426566
$ws.Range("BF66").Value = 0.86169994687064932
$ws.Range("BL66").Value = 0.95568344167515362
$ws.Range("AW67").Value = 0.71473185132578321
$ws.Range("B68").Value = 0.89389987968505802
$ws.Range("BO68").Value = 0.7673129782453173
